$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new rows before row 39 (shifts existing rows 39-52 down to 41-54)
$ws.Rows.Item(39).Resize(2).Insert()

# Populate the two newly inserted rows with the new content
# (order matters for shared-string table indexing)
$ws.Range("C38").Value = "https://hansjoerg.me/2020/02/09/tidymodels-for-machine-learning/"
$ws.Range("C39").Value = "https://www.tmwr.org/"
$ws.Range("B39").Value = "Kuhn and Silge book Tidy Modeling with R"
